# "refactor: Simplify contract handling logic in order processing script"
# Fill in the previously-blank "NV CONTRATO" (P) and "NV PEDIDO" (Q) columns
# for the two open order rows, and leave the selection on P5 as the user
# last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("P2").Value = 4600244283
$ws.Range("Q2").Value = 4503342086

$ws.Range("P3").Value = 4600244284
$ws.Range("Q3").Value = 4503342087

$null = $ws.Range("P5").Select()
